# OneTimeContractTemplate.docx edit:
#   - Turn the "FINANCIAL TERMS" table's hard-coded index / service name /
#     price cells for the 3 service rows into template placeholders.
#   - The "AEDxxx  price" cells get split into two runs: one for the
#     "[AEDn]" currency placeholder, one for the "[servicen_price]" amount
#     placeholder (matching the author's OOXML run layout exactly).
#
# The FINANCIAL TERMS table is Tables(3); rows 5/6/7 are the three
# service lines (col1=index, col2=service name, col3=AED + price).

$d = $word.ActiveDocument

function Split-RunAt($rangeStart, $rangeEnd) {
    # Toggling a character-formatting property on and back off forces the
    # engine to carve out [$rangeStart,$rangeEnd) as its own <w:r>, split
    # cleanly from its neighbours, while leaving the final formatting
    # byte-for-byte identical to what it was before.
    $r = $word.ActiveDocument.Range($rangeStart, $rangeEnd)
    $r.Font.Bold = $true
    $r.Font.Bold = $false
}

# ---------------------------------------------------------------------
# Row 5 : Service 1 (Premium AC Maintenance)
# ---------------------------------------------------------------------

$t = $d.Tables.Item(3)
$t.Cell(5, 1).Range.Find.Execute("1", $true, $true, $false, $false, $false, $true, 0, $false, "[index1]", 1) | Out-Null

$t = $d.Tables.Item(3)
$t.Cell(5, 2).Range.Find.Execute("Premium AC Maintenance ", $true, $false, $false, $false, $false, $true, 0, $false, "[service1_name]", 1) | Out-Null

$t = $d.Tables.Item(3)
$t.Cell(5, 3).Range.Find.Execute("AED   [ac_maintenance_price]", $true, $false, $false, $false, $false, $true, 0, $false, "[AED1]   [service1_price]", 1) | Out-Null

$t = $d.Tables.Item(3)
$priceCellStart = $t.Cell(5, 3).Range.Start
# Cell text is now "      [AED1]   [service1_price]" -> the spacer run is
# 6 chars, "[AED1]" is 6 chars.
$b1 = $priceCellStart + 6
$b2 = $priceCellStart + 12
Split-RunAt $b1 $b2

# ---------------------------------------------------------------------
# Row 6 : Service 2 (Premium AC Repairs)
# ---------------------------------------------------------------------

$t = $d.Tables.Item(3)
$t.Cell(6, 1).Range.Find.Execute("2", $true, $true, $false, $false, $false, $true, 0, $false, "[index2]", 1) | Out-Null

$t = $d.Tables.Item(3)
$t.Cell(6, 2).Range.Find.Execute("Premium AC Repairs", $true, $false, $false, $false, $false, $true, 0, $false, "[service2_name]", 1) | Out-Null

$t = $d.Tables.Item(3)
$t.Cell(6, 3).Range.Find.Execute("AED   [ac_repair_price]", $true, $false, $false, $false, $false, $true, 0, $false, "[AED2]   [service2_price]", 1) | Out-Null

$t = $d.Tables.Item(3)
$priceCellStart = $t.Cell(6, 3).Range.Start
# Cell text is now "[AED2]   [service2_price]" -> "[AED2]" is 6 chars; the
# run split for this row keeps a trailing space with [AED2] and only two
# leading spaces on [service2_price] ("[AED2] " + "  [service2_price]").
$b1 = $priceCellStart
$b2 = $priceCellStart + 7
Split-RunAt $b1 $b2

# ---------------------------------------------------------------------
# Row 7 : Service 3 (Other)
# ---------------------------------------------------------------------

$t = $d.Tables.Item(3)
$t.Cell(7, 1).Range.Find.Execute("3", $true, $true, $false, $false, $false, $true, 0, $false, "[index3]", 1) | Out-Null

$t = $d.Tables.Item(3)
$t.Cell(7, 2).Range.Find.Execute("Other", $true, $true, $false, $false, $false, $true, 0, $false, "[service3_name]", 1) | Out-Null

$t = $d.Tables.Item(3)
$t.Cell(7, 3).Range.Find.Execute("AED   [other_price]", $true, $false, $false, $false, $false, $true, 0, $false, "[AED3]   [service3_price]", 1) | Out-Null

$t = $d.Tables.Item(3)
$priceCellStart = $t.Cell(7, 3).Range.Start
# Cell text is now "[AED3]   [service3_price]" -> "[AED3]" is 6 chars.
$b1 = $priceCellStart
$b2 = $priceCellStart + 6
Split-RunAt $b1 $b2

Write-Host "done"
